# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 16, pushing the two
# previously-existing rows (old row 16 -> new row 17, old row 17 -> new row 18)
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 17; this shifts the former row 17
# down to row 18, and leaves the former row 16 (soon to be overwritten
# with new data) still at row 16.
$ws.Rows("17:17").Insert()

# Row 17 receives the data that used to live in row 16 before the edit
# (unchanged values, just shifted down one row).
$ws.Range("A17").Value = 5
$ws.Range("B17").Value = 'Macroferia Regional de Talca'
$ws.Range("C17").Value = 'Maule'
$ws.Range("D17").Value2 = 44447
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 'Fruta'
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = 'Otros'
$ws.Range("I17").Value = 100107002
$ws.Range("J17").Value = 'Chirimoya'
$ws.Range("K17").Value = 'Cultivar IV Región'
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 32000
$ws.Range("O17").Value = 32000
$ws.Range("P17").Value = 32000
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 3200
$ws.Range("T17").Value = 10

# Row 16 becomes the new weekly record with the updated figures.
$ws.Range("D16").Value2 = 44460
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 30000
$ws.Range("O16").Value = 30000
$ws.Range("P16").Value = 30000
$ws.Range("S16").Value = 3000
